$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update changed numeric/value cells for rows 2-9 (columns C, D, G, H)
$ws.Range("C2").Value = 99.27232944997668
$ws.Range("D2").Value = 77.06619675504865
$ws.Range("G2").Value = 0.65
$ws.Range("H2").Value = 2

$ws.Range("C3").Value = 99.31761775484055
$ws.Range("D3").Value = 77.30185823032386
$ws.Range("H3").Value = 2

$ws.Range("C4").Value = 98.67906208244496
$ws.Range("D4").Value = 75.92257468400956
$ws.Range("H4").Value = 2

$ws.Range("C5").Value = 99.29151210483801
$ws.Range("D5").Value = 77.31123566041849
$ws.Range("H5").Value = 2

$ws.Range("C6").Value = 98.82335277137378
$ws.Range("D6").Value = 77.15726338055643
$ws.Range("H6").Value = 2

$ws.Range("C7").Value = 98.7938244478818
$ws.Range("D7").Value = 77.49641283613835
$ws.Range("H7").Value = 2

$ws.Range("C8").Value = 99.09208432479556
$ws.Range("D8").Value = 77.34105247934632
$ws.Range("H8").Value = 2

$ws.Range("C9").Value = 98.56039592053354
$ws.Range("D9").Value = 75.55707046124016
$ws.Range("H9").Value = 2

# Delete row 10 entirely (was sample S9)
$ws.Rows.Item(10).Delete()
